$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing numeric metrics (B2, C2, D2)
$ws.Range("B2").Value = 0.09464907031431556
$ws.Range("C2").Value = 0.9986919242206602
$ws.Range("D2").Value = 0.2518462643775816

# Add new header "Modelo" in F1, matching the header style used by A1:E1
$ws.Range("F1").Value = "Modelo"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Add model description text in F2
$modelo = "Pipeline(steps=[('model'," + [char]10 + "                 RandomForestRegressor(max_depth=5, n_estimators=150))])"
$ws.Range("F2").Value = $modelo
